# Generate Report for Handoff
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
# - Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps bumped forward
# - Narrower "status/date" columns on Overview, zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text (was: "Handed back: in sync with en-US") ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Update timestamps ---
# Overview!G2 and de-de!H2 shared the same "Latest HO Xliff Generate Date" value
$overview.Range("G2").Value = "2016-09-02 13:09:08"
$dede.Range("H2").Value = "2016-09-02 13:09:08"

# zh-cn!H2 "Latest HO Xliff Generate Date" moved forward independently
$zhcn.Range("H2").Value = "2016-09-02 13:08:58"

# --- Narrow the status/date columns (29.9777... -> 17.2159... characters-wide) ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
